$d = $word.ActiveDocument

# Paragraph 2: both citations become (Johnson 45)
$d.Content.Find.Execute("(Ref-A1B2C3). However", $true, $false, $false, $false, $false, $true, 1, $false, "(Johnson 45). However", 2)
$d.Content.Find.Execute("2019-2020 (Ref-D4E5F6).", $true, $false, $false, $false, $false, $true, 1, $false, "2019-2020 (Johnson 45).", 2)

# Paragraph 3: (Ref-A1B2C3) -> (Ref-u278538)
$d.Content.Find.Execute("traveler (Ref-A1B2C3).", $true, $false, $false, $false, $false, $true, 1, $false, "traveler (Ref-u278538).", 2)

# Paragraph 4: both citations become (Ref-f991473)
$d.Content.Find.Execute("survey (Ref-J7Y3X2).", $true, $false, $false, $false, $false, $true, 1, $false, "survey (Ref-f991473).", 2)
$d.Content.Find.Execute("survey (Ref-Z91K0X).", $true, $false, $false, $false, $false, $true, 1, $false, "survey (Ref-f991473).", 2)

# Paragraph 5: (Ref-A1B2C3) -> (Ref-s532144)
$d.Content.Find.Execute("vehicles (Ref-A1B2C3).", $true, $false, $false, $false, $false, $true, 1, $false, "vehicles (Ref-s532144).", 2)
